# ProjectConfiguration.xlsx update:
# - Row 5 ("individualParamsFile" / "IndividualParameters.xlsx") is renamed/repurposed
#   to "individualsFile" / "Individuals.xlsx" (individual parameters now live inside
#   the "Individuals" file).
# - Row 6 ("individualPhysiologyFile" / "IndividualBiometrics.xlsx") is removed entirely,
#   since individual biometrics are now also part of the "Individuals" file.
# - All following rows shift up by one.
# - Selection moves from C6 to C5 to follow the row that used to be selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 in place (Property / Value columns; Description stays the same).
$ws.Range("A5").Value = "individualsFile"
$ws.Range("B5").Value = "Individuals.xlsx"

# Remove the old row 6 (individualPhysiologyFile / IndividualBiometrics.xlsx);
# rows below shift up automatically.
$ws.Rows("6").Delete()

# Restore the selection to match the new layout.
$ws.Range("C5").Select() | Out-Null
